# Reversed sorting in the G12C Neural Network (rows 52-62) and
# G12D Neural Network (rows 65-75) sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- G12C Neural Network block (rows 52-62) ----
$g12c = @(
    @("CHEMBL1200728", "GUANIDINE HYDROCHLORIDE", 10.190047),
    @("CHEMBL2107067", "TESTOSTERONE UNDECANOATE", 7.7600026),
    @("CHEMBL1371", "CHLORZOXAZONE", 7.4483685),
    @("CHEMBL492", "ETIDOCAINE", 5.872513),
    @("CHEMBL878", "METOLAZONE", 5.800352),
    @("CHEMBL1529", "DIPHENIDOL HYDROCHLORIDE", 5.622511),
    @("CHEMBL1200410", "PROCARBAZINE HYDROCHLORIDE", 5.469802),
    @("CHEMBL5315118", "ILOPROST TROMETHAMINE", 4.2749567),
    @("CHEMBL730", "NITROGLYCERIN", 3.3668554),
    @("CHEMBL1200968", "HYDROCORTISONE SODIUM PHOSPHATE", $null),
    @("CHEMBL1371200", "CANRENOATE POTASSIUM", $null)
)

$r = 52
foreach ($row in $g12c) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($null -ne $row[2]) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $r = $r + 1
}

# ---- G12D Neural Network block (rows 65-75) ----
$g12d = @(
    @("CHEMBL4303454", "DORIPENEM MONOHYDRATE", 6.483579),
    @("CHEMBL1200635", "HYDROCORTAMATE HYDROCHLORIDE", 6.2385273),
    @("CHEMBL2107067", "TESTOSTERONE UNDECANOATE", 5.975186),
    @("CHEMBL492", "ETIDOCAINE", 5.768978),
    @("CHEMBL878", "METOLAZONE", 5.024636),
    @("CHEMBL1529", "DIPHENIDOL HYDROCHLORIDE", 4.085992),
    @("CHEMBL5315118", "ILOPROST TROMETHAMINE", 3.448734),
    @("CHEMBL730", "NITROGLYCERIN", 3.1568189),
    @("CHEMBL1200968", "HYDROCORTISONE SODIUM PHOSPHATE", $null),
    @("CHEMBL1371200", "CANRENOATE POTASSIUM", $null),
    @("CHEMBL1200487", "ETHACRYNATE SODIUM", $null)
)

$r = 65
foreach ($row in $g12d) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($null -ne $row[2]) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    } else {
        # Force an empty Text cell (matches the already-empty cells in the
        # G12C block above) rather than clearing to a blank Number cell.
        $ws.Cells.Item($r, 3).Value = "'"
        $ws.Cells.Item($r, 3).ClearFormats()
    }
    $r = $r + 1
}
